# Update DM integration test fixture
# - Bold the header row on each sheet (CodeSchemes, Codes, Extensions)
# - Widen columns to match the new (bold-header) auto-fit widths
# - Replace a batch of UUID values with freshly generated ones

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # CodeSchemes
$ws2 = $wb.Worksheets.Item(2)   # Codes
$ws3 = $wb.Worksheets.Item(3)   # Extensions

# ---------------------------------------------------------------------------
# Column widths (character units fed to Excel's ColumnWidth property so the
# stored xlsx width lands on the desired value).
# ---------------------------------------------------------------------------

$ws1.Columns.Item(1).ColumnWidth  = 34.42857142857143
$ws1.Columns.Item(2).ColumnWidth  = 17.42857142857143
$ws1.Columns.Item(3).ColumnWidth  = 25.714285714285715
$ws1.Columns.Item(4).ColumnWidth  = 22.42857142857143
$ws1.Columns.Item(5).ColumnWidth  = 14.142857142857142
$ws1.Columns.Item(6).ColumnWidth  = 19.142857142857146
$ws1.Columns.Item(7).ColumnWidth  = 16.857142857142854
$ws1.Columns.Item(8).ColumnWidth  = 19.142857142857146
$ws1.Columns.Item(9).ColumnWidth  = 65.28571428571428
$ws1.Columns.Item(10).ColumnWidth = 24.0
$ws1.Columns.Item(11).ColumnWidth = 19.142857142857146
$ws1.Columns.Item(12).ColumnWidth = 15.857142857142858
$ws1.Columns.Item(13).ColumnWidth = 20.714285714285715
$ws1.Columns.Item(14).ColumnWidth = 27.285714285714285

$ws2.Columns.Item(1).ColumnWidth  = 37.85714285714286
$ws2.Columns.Item(2).ColumnWidth  = 17.42857142857143
$ws2.Columns.Item(3).ColumnWidth  = 15.857142857142858
$ws2.Columns.Item(4).ColumnWidth  = 14.142857142857142
$ws2.Columns.Item(5).ColumnWidth  = 15.857142857142858
$ws2.Columns.Item(6).ColumnWidth  = 19.142857142857146
$ws2.Columns.Item(7).ColumnWidth  = 20.714285714285715
$ws2.Columns.Item(8).ColumnWidth  = 24.0
$ws2.Columns.Item(9).ColumnWidth  = 19.142857142857146
$ws2.Columns.Item(10).ColumnWidth = 15.857142857142858
$ws2.Columns.Item(11).ColumnWidth = 34.42857142857143

$ws3.Columns.Item(1).ColumnWidth = 33.42857142857143
$ws3.Columns.Item(2).ColumnWidth = 17.42857142857143
$ws3.Columns.Item(3).ColumnWidth = 14.142857142857142
$ws3.Columns.Item(4).ColumnWidth = 24.0
$ws3.Columns.Item(5).ColumnWidth = 15.857142857142858
$ws3.Columns.Item(6).ColumnWidth = 19.142857142857146
$ws3.Columns.Item(7).ColumnWidth = 19.142857142857146
$ws3.Columns.Item(8).ColumnWidth = 15.857142857142858
$ws3.Columns.Item(9).ColumnWidth = 24.0

# ---------------------------------------------------------------------------
# Bold header rows
# ---------------------------------------------------------------------------

$ws1.Range("A1:N1").Font.Bold = $true
$ws2.Range("A1:K1").Font.Bold = $true
$ws3.Range("A1:I1").Font.Bold = $true

# ---------------------------------------------------------------------------
# Refresh UUID values
# ---------------------------------------------------------------------------

# CodeSchemes sheet
$ws1.Range("A2").Value = "659f916f-54b0-4e94-ad01-029d40e33885"

# Codes sheet
$ws2.Range("A2").Value = "b6034beb-48b1-476c-b389-4466a550f570"
$ws2.Range("K2").Value = "142e10e1-4e52-43f6-99e8-a21a5f8ec354"

$ws2.Range("A3").Value = "a8ab2568-627c-4f66-893f-5f1bf29a8fdb"
$ws2.Range("K3").Value = "7750dc2b-525f-4f26-a432-e4f04ad5b5e8"

$ws2.Range("A4").Value = "0d194d52-5f80-4f79-ab49-892ba58e2e57"
$ws2.Range("K4").Value = "04867dfb-9eae-4665-8d46-f66b78260752"

$ws2.Range("A5").Value = "23319de9-7bd4-45b2-860d-860585e64188"
$ws2.Range("K5").Value = "1824d625-7eec-490e-885e-e00b903f28ac"

$ws2.Range("A6").Value = "20aa51cd-f8e1-4411-abc1-52caa9ff6a6a"
$ws2.Range("K6").Value = "b45f5fe5-89c2-4561-86fd-00315aa93d3d"

$ws2.Range("A7").Value = "dccb79a2-e7fe-4fb3-afa4-7ea33f1daef0"
$ws2.Range("K7").Value = "42058b83-d1d7-4018-bd20-ea053d4f9aef"

$ws2.Range("A8").Value = "53573e19-43f9-44a9-9645-9351cb3df12a"
$ws2.Range("K8").Value = "ed5c8cc0-3299-463d-9474-56242187e817"

# Extensions sheet
$ws3.Range("A2").Value = "256de386-fd42-4816-a5b2-2840514f53ae"
